$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original headers: A1=Iccid  B1=Min  C1=Nit  D1=Mensaje
# New headers:      A1=Iccid  B1=Nit  C1=Cedula  D1=Min  E1=Mensaje
# (a "Cedula" column is inserted so it can be captured per line instead
#  of applying to the whole batch)

$min = $ws.Range("B1").Value2
$nit = $ws.Range("C1").Value2
$mensaje = $ws.Range("D1").Value2

$ws.Range("B1").Value = $nit
$ws.Range("C1").Value = "Cedula"
$ws.Range("D1").Value = $min
$ws.Range("E1").Value = $mensaje

# Match the bold/centered/bordered header style used by the other header
# cells on the newly added E1 cell
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Reflect the reported active cell/selection after the edit
$ws.Range("D6").Select()
